# =====================================================================
# Edit script for poland_iii-liga-group-i_2023-2024.xlsx
#
# The upstream scrape re-ran and:
#   1. corrected several match rows whose home/away/odds/url details had
#      been written onto the wrong physical row (rows 13-14, 22 & 24,
#      25-27, 30-31) - row 23 was already correct and is left untouched.
#   2. appended 6 new matches that were scraped afterwards (rows 68-73),
#      extending the sheet from A1:V67 to A1:V73.
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Fix mis-ordered match rows: the F:V (match detail) content of
#    several row pairs/groups was shifted to the wrong row. Restore
#    the correct content for each physical row (A-E stay as-is).
# ------------------------------------------------------------------

# Row 13 gets the match that used to be shown on row 14
$ws.Range("F13").Value = "T. Mazowiecki"
$ws.Range("H13").Value = "Legia II"
$ws.Range("K13").Value = "12/08/2023 00:12"
$ws.Range("M13").Value = "13/08/2023 11:33"
$ws.Range("O13").Value = "12/08/2023 00:12"
$ws.Range("Q13").Value = "13/08/2023 11:33"
$ws.Range("S13").Value = "12/08/2023 00:12"
$ws.Range("U13").Value = "13/08/2023 11:33"
$ws.Range("V13").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/tomaszow-mazowiecki-legia/W06rTj5f/"
$ws.Range("G13").Value = 4
$ws.Range("I13").Value = 2
$ws.Range("J13").Value = 2.72
$ws.Range("L13").Value = 3.38
$ws.Range("N13").Value = 3.27
$ws.Range("P13").Value = 3.51
$ws.Range("R13").Value = 2.16
$ws.Range("T13").Value = 1.91

# Row 14 gets the match that used to be shown on row 13
$ws.Range("F14").Value = "Pelikan"
$ws.Range("H14").Value = "Zambrow"
$ws.Range("K14").Value = "12/08/2023 00:12"
$ws.Range("M14").Value = "12/08/2023 11:16"
$ws.Range("O14").Value = "12/08/2023 00:12"
$ws.Range("Q14").Value = "13/08/2023 10:04"
$ws.Range("S14").Value = "12/08/2023 00:12"
$ws.Range("U14").Value = "12/08/2023 11:16"
$ws.Range("V14").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/pelikan-olimpia-zambrow/Ai5nSAK0/"
$ws.Range("G14").Value = 3
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 1.66
$ws.Range("L14").Value = 1.6
$ws.Range("N14").Value = 3.59
$ws.Range("P14").Value = 3.89
$ws.Range("R14").Value = 3.84
$ws.Range("T14").Value = 4.34

# Row 22 gets the match that used to be shown on row 24
$ws.Range("F22").Value = "Sulejowek"
$ws.Range("H22").Value = "Grodzisk M."
$ws.Range("K22").Value = "19/08/2023 09:27"
$ws.Range("M22").Value = "19/08/2023 16:03"
$ws.Range("O22").Value = "19/08/2023 09:27"
$ws.Range("Q22").Value = "19/08/2023 16:03"
$ws.Range("S22").Value = "19/08/2023 09:27"
$ws.Range("U22").Value = "19/08/2023 16:03"
$ws.Range("V22").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/victoria-sulejowek-gks-pogon-grodzisk-mazowiecki/IHUPuUl0/"
$ws.Range("G22").Value = 1
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 2.65
$ws.Range("L22").Value = 2.88
$ws.Range("N22").Value = 3.32
$ws.Range("P22").Value = 3.4
$ws.Range("R22").Value = 2.31
$ws.Range("T22").Value = 2.16

# Row 24 gets the match that used to be shown on row 22
$ws.Range("F24").Value = "GKS Belchatow"
$ws.Range("H24").Value = "Concordia Elblag"
$ws.Range("K24").Value = "19/08/2023 09:26"
$ws.Range("M24").Value = "19/08/2023 16:05"
$ws.Range("O24").Value = "19/08/2023 09:26"
$ws.Range("Q24").Value = "19/08/2023 16:12"
$ws.Range("S24").Value = "19/08/2023 09:26"
$ws.Range("U24").Value = "19/08/2023 16:05"
$ws.Range("V24").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/gks-belchatow-concordia-elblag/KtFO0VQP/"
$ws.Range("G24").Value = 1
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 1.46
$ws.Range("L24").Value = 1.6
$ws.Range("N24").Value = 4.19
$ws.Range("P24").Value = 3.9
$ws.Range("R24").Value = 5.27
$ws.Range("T24").Value = 4.49

# Row 25 gets the match that used to be shown on row 26
$ws.Range("F25").Value = "Skierniewice"
$ws.Range("H25").Value = "Wikielec"
$ws.Range("K25").Value = "23/08/2023 11:12"
$ws.Range("M25").Value = "23/08/2023 16:54"
$ws.Range("O25").Value = "23/08/2023 11:12"
$ws.Range("Q25").Value = "23/08/2023 16:54"
$ws.Range("S25").Value = "23/08/2023 11:12"
$ws.Range("U25").Value = "23/08/2023 16:54"
$ws.Range("V25").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/unia-skierniewice-gks-wikielec/rwu2EleP/"
$ws.Range("G25").Value = 2
$ws.Range("I25").Value = 2
$ws.Range("J25").Value = 1.66
$ws.Range("L25").Value = 1.87
$ws.Range("N25").Value = 3.79
$ws.Range("P25").Value = 3.8
$ws.Range("R25").Value = 3.9
$ws.Range("T25").Value = 3.24

# Row 26 gets the match that used to be shown on row 27
$ws.Range("F26").Value = "Swit Mazowiecki"
$ws.Range("H26").Value = "Jagiellonia II"
$ws.Range("K26").Value = "22/08/2023 05:12"
$ws.Range("M26").Value = "23/08/2023 16:11"
$ws.Range("O26").Value = "22/08/2023 05:12"
$ws.Range("Q26").Value = "23/08/2023 16:11"
$ws.Range("S26").Value = "22/08/2023 05:12"
$ws.Range("U26").Value = "23/08/2023 16:11"
$ws.Range("V26").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/swit-mazowiecki-jagiellonia/SUdFjo3I/"
$ws.Range("G26").Value = 2
$ws.Range("I26").Value = 2
$ws.Range("J26").Value = 1.74
$ws.Range("L26").Value = 1.58
$ws.Range("N26").Value = 3.58
$ws.Range("P26").Value = 4.02
$ws.Range("R26").Value = 3.42
$ws.Range("T26").Value = 4.46

# Row 27 gets the match that used to be shown on row 25
$ws.Range("F27").Value = "Concordia Elblag"
$ws.Range("H27").Value = "Zambrow"
$ws.Range("K27").Value = "22/08/2023 05:12"
$ws.Range("M27").Value = "22/08/2023 22:10"
$ws.Range("O27").Value = "22/08/2023 05:12"
$ws.Range("Q27").Value = "23/08/2023 15:00"
$ws.Range("S27").Value = "22/08/2023 05:12"
$ws.Range("U27").Value = "23/08/2023 13:49"
$ws.Range("V27").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/concordia-elblag-olimpia-zambrow/K4g7h7Y5/"
$ws.Range("G27").Value = 1
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1.97
$ws.Range("L27").Value = 1.74
$ws.Range("N27").Value = 3.27
$ws.Range("P27").Value = 3.61
$ws.Range("R27").Value = 2.98
$ws.Range("T27").Value = 3.89

# Row 30 gets the match that used to be shown on row 31
$ws.Range("F30").Value = "Grodzisk M."
$ws.Range("H30").Value = "Warta Sieradz"
$ws.Range("K30").Value = "22/08/2023 05:42"
$ws.Range("M30").Value = "23/08/2023 16:53"
$ws.Range("O30").Value = "22/08/2023 05:42"
$ws.Range("Q30").Value = "23/08/2023 16:53"
$ws.Range("S30").Value = "22/08/2023 05:42"
$ws.Range("U30").Value = "23/08/2023 16:53"
$ws.Range("V30").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/gks-pogon-grodzisk-mazowiecki-warta-sieradz/pnRT983t/"
$ws.Range("G30").Value = 3
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1.4
$ws.Range("L30").Value = 1.3
$ws.Range("N30").Value = 4.18
$ws.Range("P30").Value = 5.11
$ws.Range("R30").Value = 5
$ws.Range("T30").Value = 7.03

# Row 31 gets the match that used to be shown on row 30
$ws.Range("F31").Value = "Mlawa"
$ws.Range("H31").Value = "GKS Belchatow"
$ws.Range("K31").Value = "23/08/2023 11:12"
$ws.Range("M31").Value = "23/08/2023 16:39"
$ws.Range("O31").Value = "23/08/2023 11:12"
$ws.Range("Q31").Value = "23/08/2023 16:37"
$ws.Range("S31").Value = "23/08/2023 11:12"
$ws.Range("U31").Value = "23/08/2023 16:39"
$ws.Range("V31").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/mks-mlawa-gks-belchatow/Y7cBiRmC/"
$ws.Range("G31").Value = 3
$ws.Range("I31").Value = 4
$ws.Range("J31").Value = 2.8
$ws.Range("L31").Value = 2.8
$ws.Range("N31").Value = 3.53
$ws.Range("P31").Value = 4.11
$ws.Range("R31").Value = 2.08
$ws.Range("T31").Value = 1.98

# ------------------------------------------------------------------
# 2) Append the 6 newly-scraped matches (rows 68-73), copying the
#    number formatting used by column A (bold index) and column E
#    (date/time) from the last existing data row first.
# ------------------------------------------------------------------
$ws.Range("A67").Copy() | Out-Null
$ws.Range("A68:A73").PasteSpecial(-4122) | Out-Null
$ws.Range("E67").Copy() | Out-Null
$ws.Range("E68:E73").PasteSpecial(-4122) | Out-Null

# Row 68
$ws.Range("A68").Value = 67
$ws.Range("B68").Value = "poland"
$ws.Range("C68").Value = "iii-liga-group-i"
$ws.Range("D68").Value = "2023-2024"
$ws.Range("E68").Value = 45192.625
$ws.Range("F68").Value = "Wikielec"
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = "T. Mazowiecki"
$ws.Range("I68").Value = 1
$ws.Range("J68").Value = 2.22
$ws.Range("K68").Value = "22/09/2023 02:13"
$ws.Range("L68").Value = 2.5
$ws.Range("M68").Value = "23/09/2023 14:49"
$ws.Range("N68").Value = 3.21
$ws.Range("O68").Value = "22/09/2023 02:13"
$ws.Range("P68").Value = 3.48
$ws.Range("Q68").Value = "23/09/2023 14:49"
$ws.Range("R68").Value = 2.67
$ws.Range("S68").Value = "22/09/2023 02:13"
$ws.Range("T68").Value = 2.39
$ws.Range("U68").Value = "23/09/2023 14:49"
$ws.Range("V68").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/gks-wikielec-tomaszow-mazowiecki/nXGPNkU7/"

# Row 69
$ws.Range("A69").Value = 68
$ws.Range("B69").Value = "poland"
$ws.Range("C69").Value = "iii-liga-group-i"
$ws.Range("D69").Value = "2023-2024"
$ws.Range("E69").Value = 45192.625
$ws.Range("F69").Value = "Jagiellonia II"
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = "Sulejowek"
$ws.Range("I69").Value = 2
$ws.Range("J69").Value = 2.28
$ws.Range("K69").Value = "22/09/2023 02:13"
$ws.Range("L69").Value = 2.78
$ws.Range("M69").Value = "23/09/2023 14:46"
$ws.Range("N69").Value = 3.28
$ws.Range("O69").Value = "22/09/2023 02:13"
$ws.Range("P69").Value = 3.35
$ws.Range("Q69").Value = "23/09/2023 14:50"
$ws.Range("R69").Value = 2.48
$ws.Range("S69").Value = "22/09/2023 02:13"
$ws.Range("T69").Value = 2.15
$ws.Range("U69").Value = "23/09/2023 14:46"
$ws.Range("V69").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/jagiellonia-victoria-sulejowek/SSgCogME/"

# Row 70
$ws.Range("A70").Value = 69
$ws.Range("B70").Value = "poland"
$ws.Range("C70").Value = "iii-liga-group-i"
$ws.Range("D70").Value = "2023-2024"
$ws.Range("E70").Value = 45192.66666666666
$ws.Range("F70").Value = "Concordia Elblag"
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = "Mlawa"
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2.38
$ws.Range("K70").Value = "22/09/2023 03:13"
$ws.Range("L70").Value = 2.02
$ws.Range("M70").Value = "23/09/2023 15:43"
$ws.Range("N70").Value = 3.33
$ws.Range("O70").Value = "22/09/2023 03:13"
$ws.Range("P70").Value = 3.68
$ws.Range("Q70").Value = "23/09/2023 15:43"
$ws.Range("R70").Value = 2.35
$ws.Range("S70").Value = "22/09/2023 03:13"
$ws.Range("T70").Value = 2.96
$ws.Range("U70").Value = "23/09/2023 15:43"
$ws.Range("V70").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/concordia-elblag-mks-mlawa/Sjxaxkq1/"

# Row 71
$ws.Range("A71").Value = 70
$ws.Range("B71").Value = "poland"
$ws.Range("C71").Value = "iii-liga-group-i"
$ws.Range("D71").Value = "2023-2024"
$ws.Range("E71").Value = 45192.66666666666
$ws.Range("F71").Value = "Legionowo"
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = "Warta Sieradz"
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1.6
$ws.Range("K71").Value = "22/09/2023 03:13"
$ws.Range("L71").Value = 1.77
$ws.Range("M71").Value = "23/09/2023 15:55"
$ws.Range("N71").Value = 3.75
$ws.Range("O71").Value = "22/09/2023 03:13"
$ws.Range("P71").Value = 3.78
$ws.Range("Q71").Value = "23/09/2023 15:55"
$ws.Range("R71").Value = 3.85
$ws.Range("S71").Value = "22/09/2023 03:13"
$ws.Range("T71").Value = 3.63
$ws.Range("U71").Value = "23/09/2023 15:55"
$ws.Range("V71").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/legionowo-warta-sieradz/lxe8nZ68/"

# Row 72
$ws.Range("A72").Value = 71
$ws.Range("B72").Value = "poland"
$ws.Range("C72").Value = "iii-liga-group-i"
$ws.Range("D72").Value = "2023-2024"
$ws.Range("E72").Value = 45192.66666666666
$ws.Range("F72").Value = "Zambrow"
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = "Grodzisk M."
$ws.Range("I72").Value = 3
$ws.Range("J72").Value = 4.96
$ws.Range("K72").Value = "22/09/2023 03:13"
$ws.Range("L72").Value = 3.47
$ws.Range("M72").Value = "23/09/2023 15:42"
$ws.Range("N72").Value = 4.12
$ws.Range("O72").Value = "22/09/2023 03:13"
$ws.Range("P72").Value = 3.69
$ws.Range("Q72").Value = "23/09/2023 15:42"
$ws.Range("R72").Value = 1.43
$ws.Range("S72").Value = "22/09/2023 03:13"
$ws.Range("T72").Value = 1.83
$ws.Range("U72").Value = "23/09/2023 15:42"
$ws.Range("V72").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/olimpia-zambrow-gks-pogon-grodzisk-mazowiecki/zgcKqXjR/"

# Row 73
$ws.Range("A73").Value = 72
$ws.Range("B73").Value = "poland"
$ws.Range("C73").Value = "iii-liga-group-i"
$ws.Range("D73").Value = "2023-2024"
$ws.Range("E73").Value = 45193.5
$ws.Range("F73").Value = "Pelikan"
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = "Pilica Bialobrzegi"
$ws.Range("I73").Value = 1
$ws.Range("J73").Value = 1.57
$ws.Range("K73").Value = "22/09/2023 23:12"
$ws.Range("L73").Value = 1.28
$ws.Range("M73").Value = "24/09/2023 11:58"
$ws.Range("N73").Value = 3.76
$ws.Range("O73").Value = "22/09/2023 23:12"
$ws.Range("P73").Value = 5.17
$ws.Range("Q73").Value = "24/09/2023 11:59"
$ws.Range("R73").Value = 4.23
$ws.Range("S73").Value = "22/09/2023 23:12"
$ws.Range("T73").Value = 7.56
$ws.Range("U73").Value = "24/09/2023 11:59"
$ws.Range("V73").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/pelikan-pilica-bialobrzegi/fwvivBEl/"
